$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
# The underlying data rows were reshuffled (re-sorted) while keeping the
# other descriptive columns (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) unchanged.
$rows = @{
    2  = @{ D = 45072; M = 30; N = 15000; O = 15000; P = 15000; S = 833 }
    3  = @{ D = 45069; M = 60; N = 15000; O = 15000; P = 15000; S = 833 }
    4  = @{ D = 45096; M = 30; N = 20000; O = 20000; P = 20000; S = 1111 }
    5  = @{ D = 45061; M = 40; N = 15000; O = 15000; P = 15000; S = 833 }
    6  = @{ D = 45085; M = 30; N = 19000; O = 19000; P = 19000; S = 1056 }
    7  = @{ D = 45112; M = 20; N = 20000; O = 20000; P = 20000; S = 1111 }
    8  = @{ D = 45055; M = 50; N = 15000; O = 15000; P = 15000; S = 833 }
    9  = @{ D = 45076; M = 20; N = 15000; O = 15000; P = 15000; S = 833 }
    10 = @{ D = 45083; M = 50; N = 15000; O = 15000; P = 15000; S = 833 }
    11 = @{ D = 45092; M = 60; N = 18000; O = 19000; P = 18667; S = 1037 }
    12 = @{ D = 45084; M = 50; N = 18000; O = 19000; P = 18500; S = 1028 }
    13 = @{ D = 45111; M = 20; N = 20000; O = 20000; P = 20000; S = 1111 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("S$r").Value = $vals.S
}
